# "Generate Report for Handoff"
#
# The localization-status report is regenerated: every row that was
# previously "In Translation" has now been handed off, so its Status
# flips to "Ready for handoff" and the associated "last handoff"
# timestamps are refreshed to the moment the report was produced.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# Column that used to just fit "In Translation" now needs to fit the
# longer "Ready for handoff" label; widen it to match (mirrors Excel's
# "AutoFit Column Width" after the text grows).
$statusColWidth = 16.333333333333332

# --- Overview sheet --------------------------------------------------
# Columns "zh-cn" / "de-de" mirror each locale's current Status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-13 18:48:06"
$wsOverview.Columns("E:E").ColumnWidth = $statusColWidth
$wsOverview.Columns("F:F").ColumnWidth = $statusColWidth

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-13 18:47:56"
$wsZhCn.Columns("C:C").ColumnWidth = $statusColWidth

# --- de-de sheet ---------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-13 18:48:06"
$wsDeDe.Columns("C:C").ColumnWidth = $statusColWidth
